$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Table 1 "Film" (rows 1-5): add columns D..H (Açıklama / Süre / Yayın /
# Banner Image / Tanıtım) with per-row data.
# Header row is entered left-to-right first, then each body column is filled
# top-to-bottom before moving to the next column (matches original authoring
# order so shared-string ids line up).
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Açıklama"
$ws.Range("E1").Value = "Süre"
$ws.Range("F1").Value = "Yayın"
$ws.Range("G1").Value = "Banner Image"
$ws.Range("H1").Value = "Tanıtım"

$ws.Range("D2").Value = "açıklama1"
$ws.Range("D3").Value = "açıklama2"
$ws.Range("D4").Value = "açıklama3"
$ws.Range("D5").Value = "açıklama4"

$ws.Range("E2").Value = 120
$ws.Range("E3").Value = 111
$ws.Range("E4").Value = 214
$ws.Range("E5").Value = 119

$ws.Range("F2").Value = 2024
$ws.Range("F3").Value = 2022
$ws.Range("F4").Value = 1998
$ws.Range("F5").Value = 2004

$ws.Range("G2").Value = "1.jpg"
$ws.Range("G3").Value = "2.jpg"
$ws.Range("G4").Value = "3.jpg"
$ws.Range("G5").Value = "4.jpg"

$ws.Range("H2").Value = "1.mp4"
$ws.Range("H3").Value = "2.mp4"
$ws.Range("H4").Value = "3.mp4"
$ws.Range("H5").Value = "4.mp4"

# Header formatting (D1:H1) copied from the existing "Açıklama"-styled header
# cell D7 so no new cell styles get added to styles.xml.
$ws.Range("D7").Copy()
$ws.Range("D1:H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Table 2 "Yönetmenler" (rows 7-11): add columns E..G (Doğum T. / Cinsiyet /
# Ülke). Each column (header + body) is filled top-to-bottom before moving to
# the next column.
# ---------------------------------------------------------------------------
$ws.Range("E7").Value = "Doğum T."
$ws.Range("E8").Value = 1980
$ws.Range("E9").Value = 2000
$ws.Range("E10").Value = 1975
$ws.Range("E11").Value = 1961

$ws.Range("F7").Value = "Cinsiyet"
$ws.Range("F8").Value = "Erkek"
$ws.Range("F9").Value = "Kadın"
$ws.Range("F10").Value = "Erkek"
$ws.Range("F11").Value = "Erkek"

$ws.Range("G7").Value = "Ülke"
$ws.Range("G8").Value = "Türkiye"
$ws.Range("G9").Value = "Amerika"
$ws.Range("G10").Value = "Almanya"
$ws.Range("G11").Value = "Rusya"

$ws.Range("D7").Copy()
$ws.Range("E7:G7").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Table 3 "Oyuncular" (rows 13-17): the actor names in column C change, and
# columns E..G (Doğum T. / Cinsiyet / Ülke) are added, reusing the same
# shared strings introduced for the Yönetmenler table.
# ---------------------------------------------------------------------------
$ws.Range("E13").Value = "Doğum T."
$ws.Range("F13").Value = "Cinsiyet"
$ws.Range("G13").Value = "Ülke"

$ws.Range("C15").Value = "Cem"
$ws.Range("C14").Value = "Efe"
$ws.Range("C16").Value = "Kaan"
$ws.Range("C17").Value = "Toprak"

$ws.Range("E14").Value = 1980
$ws.Range("E15").Value = 2000
$ws.Range("E16").Value = 1975
$ws.Range("E17").Value = 1961

$ws.Range("F14").Value = "Erkek"
$ws.Range("F15").Value = "Kadın"
$ws.Range("F16").Value = "Erkek"
$ws.Range("F17").Value = "Erkek"

$ws.Range("G14").Value = "Türkiye"
$ws.Range("G15").Value = "Amerika"
$ws.Range("G16").Value = "Almanya"
$ws.Range("G17").Value = "Rusya"

$ws.Range("D13").Copy()
$ws.Range("E13:G13").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Table 4 "Kategori" (rows 19-23): brand-new table with Id / Adı / Açıklama.
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "Kategori"
$ws.Range("B19").Value = "Id"
$ws.Range("C19").Value = "Adı"
$ws.Range("D19").Value = "Açıklama"

$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "Korku"
$ws.Range("D20").Value = "açıklama1"

$ws.Range("B21").Value = 2
$ws.Range("C21").Value = "Romantik"
$ws.Range("D21").Value = "açıklama2"

$ws.Range("B22").Value = 3
$ws.Range("C22").Value = "Sevgi"
$ws.Range("D22").Value = "açıklama3"

$ws.Range("B23").Value = 4
$ws.Range("C23").Value = "Hüzün"
$ws.Range("D23").Value = "açıklama4"

# Formatting for the new table: header row styled like the other table
# headers (A/C = title/column-header style, B = id-header style,
# D = sub-header style), body rows' A column styled like the blank marker
# column used elsewhere.
$ws.Range("A13").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("B13").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("D13").Copy()
$ws.Range("D19").PasteSpecial(-4122)

$ws.Range("A14").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Column G width (new "Banner Image" column).
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 12.63

# ---------------------------------------------------------------------------
# Selection moves to M16 (matches the author's final cursor position).
# ---------------------------------------------------------------------------
[void]$ws.Range("M16").Select()
